# Updated cryptos list (price + 1h volume refresh, plus a rank swap at rows 31/32)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# "numeric" entries (plain decimals such as 567.51) must be forced to stay
# TEXT -- same as the source data -- otherwise Excel silently stores them as
# numbers and drops significant trailing zeros (e.g. "1.00" -> 1, "7.30" -> 7.3).
$updates = @(
    @{ Cell = "D2"; Value = "62.548.07"; Numeric = $false }
    @{ Cell = "E2"; Value = "  -2.81%  "; Numeric = $false }
    @{ Cell = "D3"; Value = "3.364.41"; Numeric = $false }
    @{ Cell = "E3"; Value = "  -4.16%  "; Numeric = $false }
    @{ Cell = "E4"; Value = "  +0.05%  "; Numeric = $false }
    @{ Cell = "D5"; Value = "567.51"; Numeric = $true }
    @{ Cell = "E5"; Value = "  -4.10%  "; Numeric = $false }
    @{ Cell = "D6"; Value = "124.48"; Numeric = $true }
    @{ Cell = "E6"; Value = "  -7.33%  "; Numeric = $false }
    @{ Cell = "E7"; Value = "  +0.04%  "; Numeric = $false }
    @{ Cell = "D8"; Value = "3.362.19"; Numeric = $false }
    @{ Cell = "E8"; Value = "  -4.22%  "; Numeric = $false }
    @{ Cell = "D9"; Value = "0.471"; Numeric = $true }
    @{ Cell = "E9"; Value = "  -3.47%  "; Numeric = $false }
    @{ Cell = "D10"; Value = "7.30"; Numeric = $true }
    @{ Cell = "E10"; Value = "  -4.13%  "; Numeric = $false }
    @{ Cell = "D11"; Value = "0.118"; Numeric = $true }
    @{ Cell = "E11"; Value = "  -4.85%  "; Numeric = $false }
    @{ Cell = "D12"; Value = "0.372"; Numeric = $true }
    @{ Cell = "E12"; Value = "  -4.45%  "; Numeric = $false }
    @{ Cell = "D13"; Value = "3.961.12"; Numeric = $false }
    @{ Cell = "E13"; Value = "  -3.69%  "; Numeric = $false }
    @{ Cell = "E14"; Value = "  -0.71%  "; Numeric = $false }
    @{ Cell = "D15"; Value = "3.380.48"; Numeric = $false }
    @{ Cell = "E15"; Value = "  -3.78%  "; Numeric = $false }
    @{ Cell = "D16"; Value = "0.0000170"; Numeric = $true }
    @{ Cell = "E16"; Value = "  -6.17%  "; Numeric = $false }
    @{ Cell = "D17"; Value = "62.566.53"; Numeric = $false }
    @{ Cell = "E17"; Value = "  -2.77%  "; Numeric = $false }
    @{ Cell = "D18"; Value = "24.31"; Numeric = $true }
    @{ Cell = "E18"; Value = "  -5.76%  "; Numeric = $false }
    @{ Cell = "D19"; Value = "9.37"; Numeric = $true }
    @{ Cell = "E19"; Value = "  -6.07%  "; Numeric = $false }
    @{ Cell = "D20"; Value = "5.59"; Numeric = $true }
    @{ Cell = "E20"; Value = "  -2.90%  "; Numeric = $false }
    @{ Cell = "D21"; Value = "12.98"; Numeric = $true }
    @{ Cell = "E21"; Value = "  -4.37%  "; Numeric = $false }
    @{ Cell = "D22"; Value = "370.60"; Numeric = $true }
    @{ Cell = "E22"; Value = "  -5.89%  "; Numeric = $false }
    @{ Cell = "D23"; Value = "0.550"; Numeric = $true }
    @{ Cell = "E23"; Value = "  -4.57%  "; Numeric = $false }
    @{ Cell = "D24"; Value = "3.501.50"; Numeric = $false }
    @{ Cell = "E24"; Value = "  -4.13%  "; Numeric = $false }
    @{ Cell = "D26"; Value = "71.36"; Numeric = $true }
    @{ Cell = "E26"; Value = "  -4.43%  "; Numeric = $false }
    @{ Cell = "D27"; Value = "0.0000105"; Numeric = $true }
    @{ Cell = "E27"; Value = "  -10.65%  "; Numeric = $false }
    @{ Cell = "E28"; Value = "  +0.49%  "; Numeric = $false }
    @{ Cell = "D29"; Value = "6.85"; Numeric = $true }
    @{ Cell = "E29"; Value = "  -7.48%  "; Numeric = $false }
    @{ Cell = "D30"; Value = "2.12"; Numeric = $true }
    @{ Cell = "E30"; Value = "  -6.26%  "; Numeric = $false }
    @{ Cell = "B31"; Value = "USDe"; Numeric = $false }
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; Numeric = $false }
    @{ Cell = "D31"; Value = "1.00"; Numeric = $true }
    @{ Cell = "E31"; Value = "  -0.04%  "; Numeric = $false }
    @{ Cell = "B32"; Value = "InternetComputer(DFINITY)"; Numeric = $false }
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; Numeric = $false }
    @{ Cell = "D32"; Value = "7.66"; Numeric = $true }
    @{ Cell = "E32"; Value = "  -7.52%  "; Numeric = $false }
    @{ Cell = "E33"; Value = "  -5.07%  "; Numeric = $false }
    @{ Cell = "D34"; Value = "1.37"; Numeric = $true }
    @{ Cell = "E34"; Value = "  -7.00%  "; Numeric = $false }
    @{ Cell = "D35"; Value = "3.393.24"; Numeric = $false }
    @{ Cell = "E35"; Value = "  -4.18%  "; Numeric = $false }
    @{ Cell = "D36"; Value = "22.57"; Numeric = $true }
    @{ Cell = "E36"; Value = "  -3.40%  "; Numeric = $false }
    @{ Cell = "D37"; Value = "5.15"; Numeric = $true }
    @{ Cell = "E37"; Value = "  -3.27%  "; Numeric = $false }
    @{ Cell = "D38"; Value = "164.62"; Numeric = $true }
    @{ Cell = "E38"; Value = "  -1.47%  "; Numeric = $false }
    @{ Cell = "D39"; Value = "6.57"; Numeric = $true }
    @{ Cell = "E39"; Value = "  -5.60%  "; Numeric = $false }
    @{ Cell = "D40"; Value = "1.46"; Numeric = $true }
    @{ Cell = "E40"; Value = "  -6.15%  "; Numeric = $false }
    @{ Cell = "D41"; Value = "0.0743"; Numeric = $true }
    @{ Cell = "E41"; Value = "  -5.73%  "; Numeric = $false }
    @{ Cell = "E42"; Value = "  +0.17%  "; Numeric = $false }
    @{ Cell = "D43"; Value = "0.767"; Numeric = $true }
    @{ Cell = "E43"; Value = "  -5.48%  "; Numeric = $false }
    @{ Cell = "D44"; Value = "41.26"; Numeric = $true }
    @{ Cell = "E44"; Value = "  -2.36%  "; Numeric = $false }
    @{ Cell = "D45"; Value = "4.20"; Numeric = $true }
    @{ Cell = "E45"; Value = "  -5.56%  "; Numeric = $false }
    @{ Cell = "D46"; Value = "22.56"; Numeric = $true }
    @{ Cell = "E46"; Value = "  -10.36%  "; Numeric = $false }
    @{ Cell = "D47"; Value = "1.53"; Numeric = $true }
    @{ Cell = "E47"; Value = "  -7.66%  "; Numeric = $false }
    @{ Cell = "D48"; Value = "1.06"; Numeric = $true }
    @{ Cell = "E48"; Value = "  -9.99%  "; Numeric = $false }
    @{ Cell = "D49"; Value = "6.57"; Numeric = $true }
    @{ Cell = "E49"; Value = "  -3.56%  "; Numeric = $false }
    @{ Cell = "D50"; Value = "2.235.19"; Numeric = $false }
    @{ Cell = "E50"; Value = "  -6.52%  "; Numeric = $false }
    @{ Cell = "D51"; Value = "0.842"; Numeric = $true }
    @{ Cell = "E51"; Value = "  -6.17%  "; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Force text storage so the literal digits (incl. trailing zeros) survive,
        # then restore the default style so no stray formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Host "Updated cryptos list"
